$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '42.888.97'
Set-TextValue 2 5 '  -0.42%  '

# Row 3
Set-TextValue 3 4 '2.304.27'
Set-TextValue 3 5 '  +0.07%  '

# Row 4
Set-TextValue 4 4 '0.999'
Set-TextValue 4 5 '  -0.15%  '

# Row 5
Set-TextValue 5 4 '305.95'
Set-TextValue 5 5 '  +1.91%  '

# Row 6
Set-TextValue 6 4 '96.52'
Set-TextValue 6 5 '  -1.54%  '

# Row 7
Set-TextValue 7 4 '0.508'
Set-TextValue 7 5 '  -2.06%  '

# Row 8
Set-TextValue 8 4 '0.998'
Set-TextValue 8 5 '  -0.24%  '

# Row 9
Set-TextValue 9 4 '0.503'
Set-TextValue 9 5 '  -2.57%  '

# Row 10
Set-TextValue 10 4 '35.46'
Set-TextValue 10 5 '  -1.74%  '

# Row 11
Set-TextValue 11 4 '0.0793'
Set-TextValue 11 5 '  +0.17%  '

# Row 12
Set-TextValue 12 4 '18.44'
Set-TextValue 12 5 '  +4.22%  '

# Row 13
Set-TextValue 13 4 '0.119'
Set-TextValue 13 5 '  +1.27%  '

# Row 14
Set-TextValue 14 4 '6.77'
Set-TextValue 14 5 '  -1.64%  '

# Row 15
Set-TextValue 15 4 '2.637.10'
Set-TextValue 15 5 '  -0.90%  '

# Row 16
Set-TextValue 16 4 '2.289.01'
Set-TextValue 16 5 '  +0.28%  '

# Row 17
Set-TextValue 17 4 '0.782'
Set-TextValue 17 5 '  -0.84%  '

# Row 18
Set-TextValue 18 4 '42.758.28'
Set-TextValue 18 5 '  -0.49%  '

# Row 19
Set-TextValue 19 4 '13.00'
Set-TextValue 19 5 '  +2.47%  '

# Row 20
Set-TextValue 20 4 '0.0₃0898'
Set-TextValue 20 5 '  -1.36%  '

# Row 21
Set-TextValue 21 4 '6.04'
Set-TextValue 21 5 '  -1.84%  '

# Row 22
Set-TextValue 22 4 '67.31'
Set-TextValue 22 5 '  -1.83%  '

# Row 23
Set-TextValue 23 4 '236.51'
Set-TextValue 23 5 '  -0.61%  '

# Row 24
Set-TextValue 24 5 '  -1.61%  '

# Row 25
Set-TextValue 25 2 'PancakeSwap'
Set-TextValue 25 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 25 4 '2.46'
Set-TextValue 25 5 '  +1.04%  '

# Row 26
Set-TextValue 26 2 'Dai'
Set-TextValue 26 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 26 4 '1.00'
Set-TextValue 26 5 '  +0.16%  '

# Row 27
Set-TextValue 27 5 '  +0.04%  '

# Row 28
Set-TextValue 28 4 '25.36'
Set-TextValue 28 5 '  +1.26%  '

# Row 29
Set-TextValue 29 2 'Toncoin'
Set-TextValue 29 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 29 4 '2.18'
Set-TextValue 29 5 '  +6.38%  '

# Row 30
Set-TextValue 30 2 'Monero'
Set-TextValue 30 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 30 4 '166.67'
Set-TextValue 30 5 '  +1.41%  '

# Row 31
Set-TextValue 31 4 '9.08'
Set-TextValue 31 5 '  -0.71%  '

# Row 32
Set-TextValue 32 4 '33.16'
Set-TextValue 32 5 '  +0.34%  '

# Row 33
Set-TextValue 33 5 '  +0.01%  '

# Row 34
Set-TextValue 34 4 '4.78'
Set-TextValue 34 5 '  -1.15%  '

# Row 35
Set-TextValue 35 4 '4.99'
Set-TextValue 35 5 '  -2.44%  '

# Row 36
Set-TextValue 36 4 '17.72'
Set-TextValue 36 5 '  -2.07%  '

# Row 37
Set-TextValue 37 5 '  -1.03%  '

# Row 38
Set-TextValue 38 4 '0.0694'
Set-TextValue 38 5 '  -0.38%  '

# Row 39
Set-TextValue 39 5 '  -0.49%  '

# Row 40
Set-TextValue 40 4 '1.76'
Set-TextValue 40 5 '  -1.33%  '

# Row 41
Set-TextValue 41 5 '  -0.89%  '

# Row 42
Set-TextValue 42 4 '2.71'
Set-TextValue 42 5 '  -2.56%  '

# Row 43
Set-TextValue 43 4 '2.012.79'
Set-TextValue 43 5 '  -0.41%  '

# Row 44
Set-TextValue 44 4 '0.0279'
Set-TextValue 44 5 '  -2.47%  '

# Row 45
Set-TextValue 45 5 '  +3.28%  '

# Row 46
Set-TextValue 46 4 '10.02'
Set-TextValue 46 5 '  -3.28%  '

# Row 47
Set-TextValue 47 4 '2.07'
Set-TextValue 47 5 '  -6.93%  '

# Row 48
Set-TextValue 48 4 '2.79'
Set-TextValue 48 5 '  -1.62%  '

# Row 49
Set-TextValue 49 2 'MultiversX'
Set-TextValue 49 3 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 49 4 '53.92'
Set-TextValue 49 5 '  -0.62%  '

# Row 50
Set-TextValue 50 2 'RocketPoolETH'
Set-TextValue 50 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 50 4 '2.504.28'
Set-TextValue 50 5 '  -1.03%  '

# Row 51
Set-TextValue 51 2 'BitcoinSV'
Set-TextValue 51 3 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 51 4 '71.50'
Set-TextValue 51 5 '  -2.11%  '
